$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds daily billing data for months 04/2025..06/2025
# starting on row 2 (row 1 is the header). We are adding a new block of
# data for 07/2025 at the top of the table (right after the header), which
# pushes every existing data row down by one row.
#
# To keep the existing rows' original values intact (instead of letting a
# native Insert() reformat the whole sheet), we shift the data manually,
# working from the bottom row upwards so we never overwrite a row before
# it has been read.

$lastRow = 92

for ($r = $lastRow; $r -ge 2; $r--) {
    $destRow = $r + 1
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($destRow, 4).Value = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($destRow, 5).Value = $ws.Cells.Item($r, 5).Value()
}

# New first data row: day 1 of period 07/2025
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = 17229.16
$ws.Cells.Item(2, 3).Value = 7
$ws.Cells.Item(2, 4).Value = 2025
$ws.Cells.Item(2, 5).Value = "07/2025"
